$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" (G2)
# This shared string is also referenced by de-de!H2 ("Correspond Handoff
# Datetime"), so both cells move together.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-09 06:29:14"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-11-09 06:28:59"
$wsZhCn.Range("K2").Value = "2016-11-09 06:29:52"

# de-de sheet: "Correspond Handoff Datetime" (H2) and "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-11-09 06:29:14"
$wsDeDe.Range("K2").Value = "2016-11-09 06:30:11"
